$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.097.97'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.09%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.543.54'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.15%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.05'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.16%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.99'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.42%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.13%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.07'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.74%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.03%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.57'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.24%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -4.36%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.936.12'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.18%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.511.95'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.99%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.15'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -4.02%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.853'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.47%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.130.97'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.82'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.71%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.68'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.17%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0971'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.46%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.91'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.48%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '254.08'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.62%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.97'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.09%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.06'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.48%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '27.18'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.13%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.30%  '

$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.40'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.57%  '

$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '41.13'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +4.87%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.42'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.06%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.93'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.14%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '156.58'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.75%  '

$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.38'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.75%  '

$ws.Range("B34").Value = 'Celestia'
$ws.Range("C34").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.38'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.74%  '

$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.14'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.18%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.92%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.91%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.114'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.12%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.45%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.27%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '22.01'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -6.76%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.83'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.03%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.18%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.30'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.35%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.008.13'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.16%  '

$ws.Range("B47").Value = 'BitcoinSV'
$ws.Range("C47").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '85.14'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.82%  '

$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.10'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.93%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.58'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.53%  '

$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.790.55'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.15%  '

$ws.Range("B51").Value = 'ordi'
$ws.Range("C51").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '74.93'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.73%  '
